$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for data rows 2..499 from 45192 to 45202
for ($r = 2; $r -le 499; $r++) {
    $ws.Cells.Item($r, 3).Value = 45202
}

# Row 499 picks up an explicit custom row height (matches every other data row)
$ws.Rows.Item(499).RowHeight = 15

# Append the new record as row 500
$ws.Cells.Item(500, 1).Value = "A 46189-2023"

$ws.Cells.Item(500, 2).Value = 45196
$ws.Cells.Item(500, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(500, 3).Value = 45202
$ws.Cells.Item(500, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(500, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(500, 5).Value = "VÄSTERÅS"

$ws.Cells.Item(500, 7).Value = 0.5
$ws.Cells.Item(500, 8).Value = 0
$ws.Cells.Item(500, 9).Value = 0
$ws.Cells.Item(500, 10).Value = 0
$ws.Cells.Item(500, 11).Value = 0
$ws.Cells.Item(500, 12).Value = 0
$ws.Cells.Item(500, 13).Value = 0
$ws.Cells.Item(500, 14).Value = 0
$ws.Cells.Item(500, 15).Value = 0
$ws.Cells.Item(500, 16).Value = 0
$ws.Cells.Item(500, 17).Value = 0

$ws.Cells.Item(500, 18).Value = ""
$ws.Cells.Item(500, 18).WrapText = $true
